$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: new project entry ----
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Teer Brand E-commerce Web App"
# This cell was typed/pasted without inheriting the column's border style
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "CodeSpice"
$ws.Range("D2").Value = "INT222"
$ws.Range("E2").Value = 12307911

# Plain text first (so shared-string order matches authoring order),
# hyperlinks are attached afterwards.
$ws.Range("G2").Value = "Live Link"
$ws.Range("F2").Value = "Github Link"
$ws.Range("H3").Value = "Youtube Video Presentation"
$ws.Range("H4").Value = "Client Approval Form"

# ---- Hyperlinks ----
$ws.Hyperlinks.Add($ws.Range("G2"), "https://teer-brand-mern.vercel.app/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/CodeSpice0/Teer_Brand_MERN") | Out-Null

# H2's display text typed by the user ("LinkedIn post") differs from the
# final cell text, so add the hyperlink first, then overwrite the cell text.
$ws.Hyperlinks.Add($ws.Range("H2"), "https://www.linkedin.com/posts/teer-brand-update", "", "", "LinkedIn post") | Out-Null
$ws.Range("H2").Value = "LinkedIn post (This is just post video is uploaded in youtube)"

$ws.Hyperlinks.Add($ws.Range("H3"), "https://www.youtube.com/watch?v=TeerBrandDemo") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H4"), "https://drive.google.com/drive/folders/TeerBrandClientApproval") | Out-Null

# ---- Column widths ----
$ws.Columns("F").ColumnWidth = 48.7
$ws.Columns("H").ColumnWidth = 47.9

# ---- View state ----
$ws.Range("G7").Select()
$excel.ActiveWindow.ScrollColumn = 3
